# Append a new "Machine 32" record to the machine master data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Find the first empty row right after the existing data (row 32 -> new row 33).
$lastRow = $ws.Cells($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = 10032
$ws.Cells.Item($newRow, 2).Value = "Machine 32"
$ws.Cells.Item($newRow, 3).Value = "F4-30-B9-D4-CD-6F"
$ws.Cells.Item($newRow, 4).Value = "FB5962911665"
$ws.Cells.Item($newRow, 5).Value = "192.168.0.358"
$ws.Cells.Item($newRow, 6).Value = 1001
$ws.Cells.Item($newRow, 7).Value = "eng"
$ws.Cells.Item($newRow, 8).Value = $true
$ws.Cells.Item($newRow, 9).Value = "superadmin"
$ws.Cells.Item($newRow, 10).Value = "now()"
$ws.Cells.Item($newRow, 11).Value = "now()"

# Match the author's final view state: scrolled down with C28 selected.
$ws.Range("C28").Select()
$excel.ActiveWindow.ScrollRow = 22
